# Fruta / hortaliza, semanal
# Insert the latest weekly observation as a new row right before the
# existing row 8 (i.e. between the current 6th and 7th data rows),
# pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8..46 down to 9..47, leaving a blank row 8 that inherits
# the formatting (incl. the date-cell style) of the surrounding rows.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44547
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112030
$ws.Range("G8").Value = "Poroto granado"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 22000
$ws.Range("L8").Value = 22500
$ws.Range("M8").Value = 22250
$ws.Range("N8").Value = "`$/caja 15 kilos"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 1483
$ws.Range("Q8").Value = 15
$ws.Range("R8").Value = "Hortaliza"
